# Update the "timestamp" column (Z) on the Log_Muestras sheet.
# All rows (2-29) get refreshed with the new run's timestamp values,
# matching the re-run of the pcsmote sampling log (iris D75 R75 Pentropia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2025-11-03T00:10:51.783325"
    3  = "2025-11-03T00:10:51.783325"
    4  = "2025-11-03T00:10:51.783325"
    5  = "2025-11-03T00:10:51.783325"
    6  = "2025-11-03T00:10:51.783864"
    7  = "2025-11-03T00:10:51.783864"
    8  = "2025-11-03T00:10:51.783864"
    9  = "2025-11-03T00:10:51.783864"
    10 = "2025-11-03T00:10:51.784387"
    11 = "2025-11-03T00:10:51.784387"
    12 = "2025-11-03T00:10:51.784387"
    13 = "2025-11-03T00:10:51.784387"
    14 = "2025-11-03T00:10:51.784387"
    15 = "2025-11-03T00:10:51.784387"
    16 = "2025-11-03T00:10:51.785388"
    17 = "2025-11-03T00:10:51.787711"
    18 = "2025-11-03T00:10:51.787711"
    19 = "2025-11-03T00:10:51.788271"
    20 = "2025-11-03T00:10:51.788271"
    21 = "2025-11-03T00:10:51.788271"
    22 = "2025-11-03T00:10:51.788814"
    23 = "2025-11-03T00:10:51.788814"
    24 = "2025-11-03T00:10:51.788814"
    25 = "2025-11-03T00:10:51.789357"
    26 = "2025-11-03T00:10:51.789357"
    27 = "2025-11-03T00:10:51.789357"
    28 = "2025-11-03T00:10:51.789357"
    29 = "2025-11-03T00:10:51.789357"
}

foreach ($row in $timestamps.Keys) {
    $cell = $ws.Cells.Item($row, 26)
    $cell.NumberFormat = "@"
    $cell.Value = $timestamps[$row]
}
